# dictbiz.xlsx template: add a "tenant_id" column (select list bound to
# data.findAllTenant) right before the "update_time" column, for both the
# field-comment header row (row 1) and the forRow data-model row (row 2).
#
# Before:  ... | update_usr_id | update_time
# After:   ... | update_usr_id | tenant_id | update_time

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "update_time" column is currently the last one (column K). Insert a
# new blank column in its place, which pushes the existing "update_time"
# content (and any validation/formatting on it) one column to the right,
# into column L.
$ws.Columns("K").Insert()

# Row 1: field-comment header cell driving the generated column comment +
# the data-validation dropdown for the tenant_id select list.
$ws.Range("K1").Value = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'

# Row 2: per-row data-model cell (inside the forRow loop) rendering each
# record's tenant_id label.
$ws.Range("K2").Value = '<%=model.tenant_id_lbl%>'
